$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "32.159.58"
$ws.Range("E2").Value = "  +7.68%  "
$ws.Range("D3").Value = "1.730.38"
$ws.Range("E3").Value = "  +5.44%  "
$ws.Range("D4").Value = "'0.994"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'224.74"
$ws.Range("E5").Value = "  +4.26%  "
$ws.Range("D6").Value = "'0.542"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "'30.82"
$ws.Range("E8").Value = "  +6.45%  "
$ws.Range("D9").Value = "'45.46"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("D10").Value = "'0.273"
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("D11").Value = "'0.0658"
$ws.Range("E11").Value = "  +7.96%  "
$ws.Range("D12").Value = "'0.0912"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "1.971.57"
$ws.Range("E13").Value = "  +5.14%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'10.54"
$ws.Range("E14").Value = "  +8.97%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.713.66"
$ws.Range("E15").Value = "  +4.43%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.623"
$ws.Range("E16").Value = "  +4.90%  "
$ws.Range("D17").Value = "'4.23"
$ws.Range("E17").Value = "  +8.13%  "
$ws.Range("D18").Value = "32.026.76"
$ws.Range("E18").Value = "  +7.25%  "
$ws.Range("D19").Value = "'67.96"
$ws.Range("E19").Value = "  +5.51%  "
$ws.Range("D20").Value = "'255.01"
$ws.Range("E20").Value = "  +7.16%  "
$ws.Range("D21").Value = "0.0₃0735"
$ws.Range("E21").Value = "  +4.33%  "
$ws.Range("D22").Value = "'0.996"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'10.25"
$ws.Range("E23").Value = "  +3.10%  "
$ws.Range("D24").Value = "'4.30"
$ws.Range("E24").Value = "  +4.06%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'158.95"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "'16.30"
$ws.Range("E27").Value = "  +4.34%  "
$ws.Range("E28").Value = "  +4.24%  "
$ws.Range("D29").Value = "'6.90"
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "'3.84"
$ws.Range("E31").Value = "  +13.31%  "
$ws.Range("D32").Value = "'0.0509"
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("E33").Value = "  +5.50%  "
$ws.Range("D34").Value = "'3.42"
$ws.Range("E34").Value = "  +7.13%  "
$ws.Range("D35").Value = "1.536.78"
$ws.Range("E35").Value = "  +8.14%  "
$ws.Range("E36").Value = "  +5.29%  "
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("D38").Value = "'84.49"
$ws.Range("E38").Value = "  +10.25%  "
$ws.Range("D39").Value = "'0.619"
$ws.Range("E39").Value = "  +8.40%  "
$ws.Range("E40").Value = "  +5.14%  "
$ws.Range("D41").Value = "'2.71"
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("D42").Value = "'2.32"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "'2.08"
$ws.Range("E43").Value = "  +7.15%  "
$ws.Range("D44").Value = "'0.866"
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("D45").Value = "'0.0505"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "'55.12"
$ws.Range("E46").Value = "  +8.98%  "
$ws.Range("E47").Value = "  +3.86%  "
$ws.Range("D49").Value = "1.874.59"
$ws.Range("E49").Value = "  +5.14%  "
$ws.Range("E50").Value = "  +5.57%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'94.26"
$ws.Range("E51").Value = "  +0.37%  "
